# Auto-generated edit script: updates the cryptocurrency price/volume table
# on Sheet1 of the workbook to match the refreshed data feed, and swaps two
# pairs of rows whose rank order changed (Polkadot/WrappedEther, EOS/Quant).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the inlineStr/shared-string cells
# already used throughout this sheet) even when the text looks numeric
# (e.g. "30.110.59", "1.002"), and leave the cell with no explicit style
# (ClearFormats after the write drops the transient text-format style).
function Set-TextValue($cell, [string]$value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "30.110.59"
Set-TextValue $ws.Range("E2") "  +0.39%  "
Set-TextValue $ws.Range("D3") "1.917.24"
Set-TextValue $ws.Range("E3") "  +2.60%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.07%  "
Set-TextValue $ws.Range("D5") "320.08"
Set-TextValue $ws.Range("E5") "  +0.28%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  +0.04%  "
Set-TextValue $ws.Range("D7") "0.5063"
Set-TextValue $ws.Range("E7") "  -0.62%  "
Set-TextValue $ws.Range("D8") "0.4076"
Set-TextValue $ws.Range("E8") "  +3.86%  "
Set-TextValue $ws.Range("D9") "0.08351"
Set-TextValue $ws.Range("E9") "  +2.14%  "
Set-TextValue $ws.Range("D10") "42.39"
Set-TextValue $ws.Range("E10") "  +0.71%  "
Set-TextValue $ws.Range("D11") "1.109"
Set-TextValue $ws.Range("E11") "  +1.73%  "
Set-TextValue $ws.Range("D12") "24.03"
Set-TextValue $ws.Range("E12") "  +5.77%  "
Set-TextValue $ws.Range("B13") "WrappedEther"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D13") "1.918.69"
Set-TextValue $ws.Range("E13") "  +2.79%  "
Set-TextValue $ws.Range("B14") "Polkadot"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D14") "6.414"
Set-TextValue $ws.Range("E14") "  +2.70%  "
Set-TextValue $ws.Range("D15") "7.246"
Set-TextValue $ws.Range("E15") "  +1.39%  "
Set-TextValue $ws.Range("D16") "1.003"
Set-TextValue $ws.Range("E16") "  +0.17%  "
Set-TextValue $ws.Range("D17") "92.55"
Set-TextValue $ws.Range("E17") "  +1.03%  "
Set-TextValue $ws.Range("D18") "0.00001096"
Set-TextValue $ws.Range("E18") "  +1.64%  "
Set-TextValue $ws.Range("D19") "0.06513"
Set-TextValue $ws.Range("E19") "  +2.26%  "
Set-TextValue $ws.Range("D20") "18.55"
Set-TextValue $ws.Range("E20") "  +3.93%  "
Set-TextValue $ws.Range("E21") "  -0.01%  "
Set-TextValue $ws.Range("D22") "5.944"
Set-TextValue $ws.Range("E22") "  +2.73%  "
Set-TextValue $ws.Range("D23") "30.123.81"
Set-TextValue $ws.Range("E23") "  +0.51%  "
Set-TextValue $ws.Range("E24") "  +2.65%  "
Set-TextValue $ws.Range("E25") "  +1.12%  "
Set-TextValue $ws.Range("D26") "2.135.48"
Set-TextValue $ws.Range("E26") "  +2.57%  "
Set-TextValue $ws.Range("D27") "21.82"
Set-TextValue $ws.Range("E27") "  +4.43%  "
Set-TextValue $ws.Range("D28") "162.85"
Set-TextValue $ws.Range("E28") "  +1.12%  "
Set-TextValue $ws.Range("D29") "2.273"
Set-TextValue $ws.Range("E29") "  +2.63%  "
Set-TextValue $ws.Range("D30") "128.88"
Set-TextValue $ws.Range("E30") "  +1.63%  "
Set-TextValue $ws.Range("D31") "1.143"
Set-TextValue $ws.Range("E31") "  +9.24%  "
Set-TextValue $ws.Range("E32") "  +1.20%  "
Set-TextValue $ws.Range("D33") "5.948"
Set-TextValue $ws.Range("E33") "  +0.92%  "
Set-TextValue $ws.Range("D34") "3.787"
Set-TextValue $ws.Range("E34") "  +1.47%  "
Set-TextValue $ws.Range("D35") "0.02456"
Set-TextValue $ws.Range("E35") "  +1.94%  "
Set-TextValue $ws.Range("D36") "5.364"
Set-TextValue $ws.Range("E36") "  +3.17%  "
Set-TextValue $ws.Range("D37") "0.06449"
Set-TextValue $ws.Range("E37") "  +2.17%  "
Set-TextValue $ws.Range("E38") "  +1.06%  "
Set-TextValue $ws.Range("D39") "0.6522"
Set-TextValue $ws.Range("E39") "  +3.84%  "
Set-TextValue $ws.Range("D40") "1.197"
Set-TextValue $ws.Range("E40") "  +2.11%  "
Set-TextValue $ws.Range("D41") "8.622"
Set-TextValue $ws.Range("E41") "  +1.76%  "
Set-TextValue $ws.Range("E42") "  +1.51%  "
Set-TextValue $ws.Range("D43") "1.210"
Set-TextValue $ws.Range("E43") "  +0.51%  "
Set-TextValue $ws.Range("D44") "13.36"
Set-TextValue $ws.Range("E44") "  +4.25%  "
Set-TextValue $ws.Range("D45") "0.6078"
Set-TextValue $ws.Range("E45") "  +3.46%  "
Set-TextValue $ws.Range("D46") "2.190"
Set-TextValue $ws.Range("E46") "  +10.55%  "
Set-TextValue $ws.Range("D47") "3.623"
Set-TextValue $ws.Range("E47") "  -0.02%  "
Set-TextValue $ws.Range("B48") "Quant"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D48") "122.31"
Set-TextValue $ws.Range("E48") "  -0.11%  "
Set-TextValue $ws.Range("B49") "EOS"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue $ws.Range("D49") "1.210"
Set-TextValue $ws.Range("E49") "  +0.79%  "
Set-TextValue $ws.Range("D50") "79.08"
Set-TextValue $ws.Range("E50") "  +3.58%  "
Set-TextValue $ws.Range("E51") "  -0.39%  "
